# Update cryptocurrency price/volume figures (Price in column D, Volume(1h)
# in column E) on the active sheet, matching the latest scrape snapshot.
#
# Values are written as plain text (matching the original cell contents,
# which are text, not numbers — e.g. "28.022.25" / "1.862.77" use dots as
# thousands separators, and the rows mix locales). Any replacement text
# that *looks* like a plain decimal number (e.g. "312.14") is written with
# a leading apostrophe so Excel stores it verbatim as text instead of
# silently parsing/rounding it as a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "28.022.25"),
    @("E2", "  +0.26%  "),
    @("D3", "1.862.77"),
    @("E4", "  +0.25%  "),
    @("D5", "312.14"),
    @("E5", "  -0.07%  "),
    @("E6", "  +0.17%  "),
    @("D7", "0.5095"),
    @("E7", "  +1.35%  "),
    @("D8", "0.3835"),
    @("E8", "  +0.39%  "),
    @("D9", "0.08305"),
    @("D10", "1.114"),
    @("E10", "  -0.26%  "),
    @("E11", "  -0.12%  "),
    @("D12", "6.225"),
    @("E12", "  -2.16%  "),
    @("D13", "20.59"),
    @("E13", "  -0.35%  "),
    @("D14", "1.857.97"),
    @("E14", "  -0.91%  "),
    @("D15", "7.217"),
    @("E15", "  -0.16%  "),
    @("D16", "1.003"),
    @("E17", "  -0.09%  "),
    @("D18", "90.93"),
    @("E18", "  -0.09%  "),
    @("D19", "0.06630"),
    @("E19", "  -0.28%  "),
    @("D20", "17.72"),
    @("E20", "  -2.29%  "),
    @("E21", "  +0.16%  "),
    @("D22", "6.042"),
    @("E22", "  -1.10%  "),
    @("D23", "28.053.08"),
    @("E23", "  +0.28%  "),
    @("E24", "  -3.45%  "),
    @("D25", "2.231"),
    @("E25", "  -1.80%  "),
    @("E26", "  +2.41%  "),
    @("D27", "2.072.64"),
    @("E27", "  -0.76%  "),
    @("D28", "158.06"),
    @("E28", "  +0.21%  "),
    @("D29", "20.55"),
    @("E29", "  -0.63%  "),
    @("D30", "124.70"),
    @("E30", "  -1.11%  "),
    @("D31", "0.1054"),
    @("E31", "  -0.88%  "),
    @("D32", "1.040"),
    @("E32", "  -1.16%  "),
    @("D33", "5.891"),
    @("E33", "  +5.37%  "),
    @("D34", "3.599"),
    @("E34", "  -0.05%  "),
    @("D35", "9.469"),
    @("E35", "  +0.28%  "),
    @("D36", "0.02420"),
    @("E36", "  +1.00%  "),
    @("D37", "0.06528"),
    @("E37", "  -0.71%  "),
    @("D38", "0.2175"),
    @("E38", "  -0.18%  "),
    @("D39", "1.207"),
    @("E39", "  +0.41%  "),
    @("D40", "0.6473"),
    @("E40", "  +1.71%  "),
    @("E41", "  -4.28%  "),
    @("D42", "4.949"),
    @("E42", "  +1.31%  "),
    @("E43", "  -2.25%  "),
    @("D44", "0.6102"),
    @("E44", "  +1.78%  "),
    @("D45", "13.09"),
    @("E45", "  -0.31%  "),
    @("D46", "1.287"),
    @("E46", "  +0.27%  "),
    @("D47", "3.670"),
    @("D48", "2.020"),
    @("E48", "  +1.48%  "),
    @("D49", "1.207"),
    @("E49", "  -1.95%  "),
    @("D50", "120.20"),
    @("E50", "  -0.31%  "),
    @("D51", "78.22"),
    @("E51", "  -0.98%  ")
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newValue = $u[1]

    if ($newValue -match '^[+-]?\d+(\.\d+)?$') {
        # Plain-looking decimal number: force text storage so Excel does
        # not reinterpret/round it as a numeric value.
        $ws.Range($cellRef).Value = "'" + $newValue
    } else {
        $ws.Range($cellRef).Value = $newValue
    }
}
